$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "66.052.68"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "3.487.49"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "604.16"
$ws.Range("E5").Value = "  +0.67%  "
Set-TextValue "D6" "143.26"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("D7").Value = "3.486.11"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.12%  "
Set-TextValue "D9" "0.475"
$ws.Range("E9").Value = "  -0.68%  "
Set-TextValue "D10" "8.19"
$ws.Range("E10").Value = "  +5.78%  "
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").Value = "4.076.47"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  -4.29%  "
Set-TextValue "D15" "30.38"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "3.484.96"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.144.83"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.117"
$ws.Range("E18").Value = "  -0.34%  "
Set-TextValue "D19" "10.38"
$ws.Range("E19").Value = "  +2.36%  "
Set-TextValue "D20" "6.16"
$ws.Range("E20").Value = "  -3.47%  "
Set-TextValue "D21" "14.76"
$ws.Range("E21").Value = "  -2.70%  "
Set-TextValue "D22" "421.08"
$ws.Range("E22").Value = "  -2.55%  "
Set-TextValue "D23" "0.591"
$ws.Range("E23").Value = "  -2.35%  "
Set-TextValue "D24" "77.49"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "3.618.29"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -4.15%  "
Set-TextValue "D28" "9.33"
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("E30").Value = "  -0.68%  "
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  -0.14%  "
Set-TextValue "D32" "0.163"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("E33").Value = "  -6.94%  "
Set-TextValue "D34" "25.20"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "3.483.13"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("E38").Value = "  -5.57%  "
Set-TextValue "D39" "7.68"
$ws.Range("E39").Value = "  -2.64%  "
Set-TextValue "D40" "0.999"
$ws.Range("E40").Value = "  -0.03%  "
Set-TextValue "D41" "170.45"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("E44").Value = "  -5.54%  "
$ws.Range("E45").Value = "  -7.17%  "
Set-TextValue "D46" "45.08"
$ws.Range("E46").Value = "  -2.59%  "
Set-TextValue "D47" "26.09"
$ws.Range("E47").Value = "  -10.04%  "
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  -4.33%  "
Set-TextValue "D51" "0.931"
